$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 6.25
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 1.6
$ws.Range("Z2").Value = 67

# Row 3 updates
$ws.Range("G3").Value = 3.4
$ws.Range("H3").Value = 2.9
$ws.Range("J3").Value = 4.33
$ws.Range("K3").Value = 1.83
$ws.Range("L3").Value = 3.25
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 7.5
$ws.Range("AB3").Value = 51
$ws.Range("AC3").Value = 6
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AI3").Value = 9.5
$ws.Range("AQ3").Value = 81
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 9.5
$ws.Range("AW3").Value = 4
